# fix(excel): properly handling ion number format retrieved as Long, BigDecimal or BigIntegers
#
# Column A (policyID) gets re-typed from numeric to text (values like
# "119736" can't safely round-trip through a double once they come from
# Ion as Long/BigInteger), and the row-2 sample that demonstrates a
# BigInteger-sized value (eq_site_limit) is split into:
#   D2 -> the in-range portion, kept numeric, shown in scientific notation
#   E2 -> the full, too-big-for-a-double literal, kept as text
# A small monospace font (JetBrains Mono, teal) highlights the reformatted
# numeric-ish text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- number formats -------------------------------------------------
$textFormat = "@"
$sciFormat  = "0.00E+00"
$dateFormat = "mm/dd/yy"

# --- column A: policyID header + values become text -----------------
$ws.Range("A1:A6").NumberFormat = $textFormat

$ws.Range("A2").Value = "119736"
$ws.Range("A3").Value = "448094"
$ws.Range("A4").Value = "206893"
$ws.Range("A5").Value = "333743"
$ws.Range("A6").Value = "172534"

# --- row 2: eq_site_limit overflowed a BigInteger --------------------
# Distinguish the two reformatted cells with a small monospace teal font
# first, so the later NumberFormat writes don't fork extra unused styles.
$ws.Range("D2:E2").Font.Name = "JetBrains Mono"
$ws.Range("D2:E2").Font.Size = 10
$ws.Range("D2:E2").Font.Family = 3
$ws.Range("D2:E2").Font.Color = 12102698  # RGB(0x2A, 0xAC, 0xB8) as BGR OLE color

# D2 keeps the truncated/representable numeric value, in scientific
# notation; E2 carries the full original (too big for IEEE-754) as text.
$ws.Range("D2").NumberFormat = $sciFormat
$ws.Range("D2").Value = 922337203685478

$ws.Range("E2").NumberFormat = $textFormat
$ws.Range("E2").Value = "9223372036854775808"

# --- cosmetic layout changes --------------------------------------
# (ColumnWidth inputs are nudged so the engine's internal character-grid
# rounding lands on the value closest to the target width.)
$ws.Columns.Item(1).ColumnWidth = 22.33
$ws.Columns.Item(2).ColumnWidth = 9.0
$ws.Columns.Item(4).ColumnWidth = 22.0
$ws.Columns.Item(6).ColumnWidth = 10.33

$ws.Rows.Item(1).RowHeight = 13.8

$ws.Range("D9").Select()
